$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 1 - headers (rename / add columns)
# ------------------------------------------------------------------
$ws.Range("D1").Value = "Delay Counter"
$ws.Range("E1").Value = "Chosen Delay"
$ws.Range("F1").Value = "Scaled Loop"
$ws.Range("G1").Value = "Chosen Scale"
$ws.Range("H1").Value = "Time shaved per register digit"
$ws.Range("I1").Value = "Register Initial value"

# ------------------------------------------------------------------
# Row 2
# ------------------------------------------------------------------
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4

$ws.Range("F2").NumberFormat = "0.00E+00"
$ws.Range("F2").Formula = '=$A$6*C2'

$ws.Range("G2").NumberFormat = "General"
$ws.Range("G2").Formula = '=C4'

$ws.Range("H2").NumberFormat = "0.00E+00"
$ws.Range("H2").Formula = '=(($A$13/G2)/256)*E2'

$ws.Range("I2").NumberFormat = "0.00"
$ws.Range("I2").Formula = '=(A13-A15)/H2'

# ------------------------------------------------------------------
# Row 3
# ------------------------------------------------------------------
$ws.Range("D3").Value = 2
$ws.Range("E3").ClearContents()
$ws.Range("E3").NumberFormat = "0.00"
$ws.Range("F3").NumberFormat = "0.00E+00"
$ws.Range("F3").Formula = '=$A$6*C3'
$ws.Range("G3").ClearContents()
$ws.Range("G3").NumberFormat = "General"
$ws.Range("H3").ClearContents()
$ws.Range("H3").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 4
# ------------------------------------------------------------------
$ws.Range("D4").Value = 3
$ws.Range("E4").ClearContents()
$ws.Range("E4").NumberFormat = "0.00"
$ws.Range("F4").NumberFormat = "0.00E+00"
$ws.Range("F4").Formula = '=$A$6*C4'
$ws.Range("G4").ClearContents()
$ws.Range("G4").NumberFormat = "General"
$ws.Range("H4").ClearContents()
$ws.Range("H4").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 5
# ------------------------------------------------------------------
$ws.Range("D5").Value = 4
$ws.Range("E5").ClearContents()
$ws.Range("E5").NumberFormat = "0.00"
$ws.Range("F5").NumberFormat = "0.00E+00"
$ws.Range("F5").Formula = '=$A$6*C5'
$ws.Range("G5").ClearContents()
$ws.Range("G5").NumberFormat = "General"
$ws.Range("H5").ClearContents()
$ws.Range("H5").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 6
# ------------------------------------------------------------------
$ws.Range("D6").Value = 5
$ws.Range("E6").ClearContents()
$ws.Range("E6").NumberFormat = "0.00"
$ws.Range("F6").NumberFormat = "0.00E+00"
$ws.Range("F6").Formula = '=$A$6*C6'
$ws.Range("G6").ClearContents()
$ws.Range("G6").NumberFormat = "General"
$ws.Range("H6").ClearContents()
$ws.Range("H6").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 7
# ------------------------------------------------------------------
$ws.Range("D7").Value = 6
$ws.Range("E7").ClearContents()
$ws.Range("E7").NumberFormat = "0.00"
$ws.Range("F7").NumberFormat = "0.00E+00"
$ws.Range("F7").Formula = '=$A$6*C7'
$ws.Range("G7").ClearContents()
$ws.Range("G7").NumberFormat = "General"
$ws.Range("H7").ClearContents()
$ws.Range("H7").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 8
# ------------------------------------------------------------------
$ws.Range("D8").Value = 7
$ws.Range("E8").ClearContents()
$ws.Range("E8").NumberFormat = "0.00"
$ws.Range("F8").NumberFormat = "0.00E+00"
$ws.Range("F8").Formula = '=$A$6*C8'
$ws.Range("G8").ClearContents()
$ws.Range("G8").NumberFormat = "General"
$ws.Range("H8").ClearContents()
$ws.Range("H8").NumberFormat = "0.00E+00"

# ------------------------------------------------------------------
# Row 9  (note: no H9 cell in the result)
# ------------------------------------------------------------------
$ws.Range("D9").Value = 8
$ws.Range("E9").ClearContents()
$ws.Range("E9").NumberFormat = "0.00"
$ws.Range("F9").NumberFormat = "0.00E+00"
$ws.Range("F9").Formula = '=$A$6*C9'
$ws.Range("G9").ClearContents()
$ws.Range("G9").NumberFormat = "General"

# ------------------------------------------------------------------
# Row 10
# ------------------------------------------------------------------
$ws.Range("D10").Value = 9

# ------------------------------------------------------------------
# Row 11 (new)
# ------------------------------------------------------------------
$ws.Range("D11").Value = 10

# ------------------------------------------------------------------
# Row 12 (new) - "Measured Loop" label
# ------------------------------------------------------------------
$ws.Range("A12").Value = "Measured Loop"
$ws.Range("D12").Value = 11

# ------------------------------------------------------------------
# Row 13 (new) - measured loop value
# ------------------------------------------------------------------
$ws.Range("A13").NumberFormat = "0.00E+00"
$ws.Range("A13").Value = 0.0017
$ws.Range("D13").Value = 12

# ------------------------------------------------------------------
# Row 14 (new) - "Required Loop" label
# ------------------------------------------------------------------
$ws.Range("A14").Value = "Required Loop"
$ws.Range("D14").Value = 13

# ------------------------------------------------------------------
# Row 15 (new) - required loop value
# ------------------------------------------------------------------
$ws.Range("A15").NumberFormat = "0.00E+00"
$ws.Range("A15").Value = 0.001
$ws.Range("D15").Value = 14

# ------------------------------------------------------------------
# Row 16 (new)
# ------------------------------------------------------------------
$ws.Range("D16").Value = 15

# ------------------------------------------------------------------
# Column widths for newly-visible / resized columns
# ------------------------------------------------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 13.85546875
$ws.Range("G1").EntireColumn.ColumnWidth = 16
$ws.Range("H1").EntireColumn.ColumnWidth = 29.140625
$ws.Range("I1").EntireColumn.ColumnWidth = 19.140625
$ws.Range("J1").EntireColumn.ColumnWidth = 16.28515625
$ws.Range("K1").EntireColumn.ColumnWidth = 10.85546875
$ws.Range("M1").EntireColumn.ColumnWidth = 11.5703125

# ------------------------------------------------------------------
# Selection moves to E3
# ------------------------------------------------------------------
$ws.Range("E3").Select()
